$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("testConstructor")
$ws2 = $wb.Worksheets.Item("testGetterSetter")
$ws2.Activate()

# ===== Row 1: header =====
$ws2.Cells.Item(1,1).Value = "function name"
$ws2.Cells.Item(1,2).Value = "id"
$ws2.Cells.Item(1,3).Value = "params"
$ws2.Cells.Item(1,5).Value = "Expected Results"

# ===== setTitle group (rows 2-4): fill A,B,C first, then E =====
$ws2.Cells.Item(2,1).Value = "setTitle"
$ws2.Cells.Item(2,2).Value = 0
$ws2.Cells.Item(2,3).Value = """"""
$ws2.Cells.Item(3,2).Value = 1
$ws2.Cells.Item(3,3).Value = "null"
$ws2.Cells.Item(4,2).Value = 2
$ws2.Cells.Item(4,3).Value = """valid title"""
$ws2.Cells.Item(2,5).Value = "IllegalArgumentException"
$ws2.Cells.Item(3,5).Value = "IllegalArgumentException"
$ws2.Cells.Item(4,5).Value = "No Action"

# ===== setPreparationTime group (rows 5-7) =====
$ws2.Cells.Item(5,1).Value = "setPreparationTime"
$ws2.Cells.Item(5,2).Value = 0
$ws2.Cells.Item(5,3).Value = 1
$ws2.Cells.Item(5,5).Value = "No Action"
$ws2.Cells.Item(6,2).Value = 1
$ws2.Cells.Item(6,3).Value = 0
$ws2.Cells.Item(6,5).Value = "IllegalArgumentException"
$ws2.Cells.Item(7,2).Value = 2
$ws2.Cells.Item(7,3).Value = -1
$ws2.Cells.Item(7,5).Value = "IllegalArgumentException"

# ===== setServings group (rows 8-10) =====
$ws2.Cells.Item(8,1).Value = "setServings"
$ws2.Cells.Item(8,2).Value = 0
$ws2.Cells.Item(8,3).Value = 1
$ws2.Cells.Item(8,5).Value = "No Action"
$ws2.Cells.Item(9,2).Value = 1
$ws2.Cells.Item(9,3).Value = 0
$ws2.Cells.Item(9,5).Value = "IllegalArgumentException"
$ws2.Cells.Item(10,2).Value = 2
$ws2.Cells.Item(10,3).Value = -1
$ws2.Cells.Item(10,5).Value = "IllegalArgumentException"

# ===== setCategory group (rows 11-13) =====
$ws2.Cells.Item(11,1).Value = "setCategory"
$ws2.Cells.Item(11,2).Value = 0
$ws2.Cells.Item(11,3).Value = """valid cat"""
$ws2.Cells.Item(11,5).Value = "No Action"
$ws2.Cells.Item(12,2).Value = 1
$ws2.Cells.Item(12,3).Value = """"""
$ws2.Cells.Item(12,5).Value = "IllegalArgumentException"
$ws2.Cells.Item(13,2).Value = 2
$ws2.Cells.Item(13,3).Value = "null"
$ws2.Cells.Item(13,5).Value = "IllegalArgumentException"

# ===== Fill formatting (copy existing fills from sheet1 rows, reuses style indices 2-5) =====
$ws1.Range("A3").Copy()
$ws2.Range("A2:F4").PasteSpecial(-4122)

$ws1.Range("A5").Copy()
$ws2.Range("A5:F7").PasteSpecial(-4122)

$ws1.Range("A7").Copy()
$ws2.Range("A8:F10").PasteSpecial(-4122)

$ws1.Range("A9").Copy()
$ws2.Range("A11:F13").PasteSpecial(-4122)

# ===== Header row bold (creates new font + style) =====
$ws2.Range("A1:F1").Font.Bold = $true

# ===== Column width =====
$ws2.Columns.Item(1).ColumnWidth = 23

# ===== Selection / active cell =====
$ws2.Range("Q27").Select()
